$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 412, shifting the existing
# rows 412-436 down to 414-438 (matches the target diff's row insertion).
$ws.Rows("412:413").Insert()

# Row 412 (new): Escarola / Primera, Fecha 2022-01-24 (serial 44585)
$ws.Range("A412").Value = 4
$ws.Range("B412").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C412").Value = "Los Lagos"
$ws.Range("D412").Value = 44585
$ws.Range("E412").Value = 10
$ws.Range("F412").Value = 100112033
$ws.Range("G412").Value = "Lechuga"
$ws.Range("H412").Value = "Escarola"
$ws.Range("I412").Value = "Primera"
$ws.Range("J412").Value = 100
$ws.Range("K412").Value = 13000
$ws.Range("L412").Value = 13000
$ws.Range("M412").Value = 13000
$ws.Range("N412").Value = "`$/caja 15 unidades"
$ws.Range("O412").Value = "Región Metropolitana"
$ws.Range("P412").Value = 867
$ws.Range("Q412").Value = 15
$ws.Range("R412").Value = "Hortaliza"

# Row 413 (new): Escarola / Segunda, Fecha 2022-01-24 (serial 44585)
$ws.Range("A413").Value = 4
$ws.Range("B413").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C413").Value = "Los Lagos"
$ws.Range("D413").Value = 44585
$ws.Range("E413").Value = 10
$ws.Range("F413").Value = 100112033
$ws.Range("G413").Value = "Lechuga"
$ws.Range("H413").Value = "Escarola"
$ws.Range("I413").Value = "Segunda"
$ws.Range("J413").Value = 100
$ws.Range("K413").Value = 10000
$ws.Range("L413").Value = 10000
$ws.Range("M413").Value = 10000
$ws.Range("N413").Value = "`$/caja 18 unidades"
$ws.Range("O413").Value = "Región Metropolitana"
$ws.Range("P413").Value = 556
$ws.Range("Q413").Value = 18
$ws.Range("R413").Value = "Hortaliza"

# Ensure the date cells keep the same custom date number format used
# throughout column D (numFmt "YYYY-MM-DD HH:MM:SS").
$ws.Range("D412:D413").NumberFormat = $ws.Range("D414").NumberFormat
